$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B,C,D,E,F,I,J,K,L,M,N correspond to Excel column indices below (G/H unchanged)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @{}
$data[2] = @(1.02, 1.053847231077304, 1.058805806406875, 1.066879532602205, 1.072929736507548, 1.048817468724508, 1.058862427238636, 1.061537036081128, 1.069588895151856, 1.075622945727848, 1.060366133338532)
$data[3] = @(1.02, 1.054742590320598, 1.059498739917605, 1.067723206866749, 1.073766675267965, 1.049042978546626, 1.059409094889193, 1.062044381979079, 1.070248177186204, 1.076276672889828, 1.060913577319857)
$data[4] = @(1.02, 1.05532267906609, 1.059947745310006, 1.068270156656193, 1.074309229847248, 1.049188112050635, 1.059762883915934, 1.062372619905345, 1.070675162312738, 1.076700024608499, 1.061267868767596)
$data[5] = @(1.02, 1.055566721972093, 1.060136656572038, 1.068500340605976, 1.074537557184756, 1.049248937207424, 1.059911629580134, 1.062510598199828, 1.070854757690629, 1.07687808312971, 1.061416825667701)
$data[6] = @(1.02, 1.055607707985835, 1.060168384315763, 1.068539003892623, 1.074575908217117, 1.049259138917536, 1.059936605343508, 1.062533764568139, 1.070884917828352, 1.07690798467206, 1.061441836899524)
$data[7] = @(1.02, 1.055325939299109, 1.059950268966631, 1.068273231420016, 1.074312279837808, 1.049188925542738, 1.059764871414213, 1.062374463630436, 1.070677561719911, 1.076702403516979, 1.061269859088351)
$data[8] = @(1.02, 1.054149669564929, 1.059039854849246, 1.067164440655485, 1.073212375882892, 1.048893843241581, 1.059047163266421, 1.061708505355262, 1.069811621881344, 1.075843803228139, 1.060551131712667)
$data[9] = @(1.02, 1.052082600977585, 1.057440496687283, 1.065218624443756, 1.071281938374146, 1.04836787823633, 1.057782973110943, 1.060534683308834, 1.068288743465133, 1.074333560284504, 1.059285146262163)
$data[10] = @(1.02, 1.050708456716209, 1.056377658218578, 1.063926909214781, 1.070000288786057, 1.048013247566, 1.056940592457271, 1.059751995579499, 1.06727560837766, 1.073328650298563, 1.058441569331341)
$data[11] = @(1.02, 1.050114381312047, 1.055918265164034, 1.063368908281199, 1.069446602524543, 1.047858751566771, 1.056575945896657, 1.059413065696093, 1.066837430389362, 1.072893987131827, 1.058076404930785)
$data[12] = @(1.02, 1.04989385770303, 1.055747751494567, 1.063161842069549, 1.06924113242857, 1.047801224551606, 1.056440517680973, 1.05928717015829, 1.066674750615522, 1.072732605931964, 1.057940784391501)
$data[13] = @(1.02, 1.049941154295507, 1.055784321554884, 1.063206249387014, 1.069285197669412, 1.047813570630425, 1.056469566669464, 1.059314175243772, 1.06670964237855, 1.072767219462187, 1.057969874632888)
$data[14] = @(1.02, 1.050096149847133, 1.055904167883581, 1.063351788029394, 1.069429614341473, 1.047853999224287, 1.056564750988694, 1.059402659162914, 1.06682398160748, 1.072880645832846, 1.058065194124769)
$data[15] = @(1.02, 1.050191666642719, 1.055978025825596, 1.06344148577067, 1.069518619933556, 1.047878890041884, 1.056623399579193, 1.05945717681323, 1.066894440298348, 1.072950541156594, 1.058123926002993)
$data[16] = @(1.02, 1.050747903437804, 1.056408164109178, 1.063963969871415, 1.070037062172139, 1.0480234812406, 1.056964795269322, 1.059774488905839, 1.067304699810254, 1.073357507482566, 1.058465806514161)
$data[17] = @(1.02, 1.051097068307654, 1.056678200047138, 1.064292065107743, 1.070362610350544, 1.048113928680681, 1.057178974052744, 1.059973525680849, 1.06756218401274, 1.073612913473777, 1.058680289456023)
$data[18] = @(1.02, 1.051300820683532, 1.056835786735132, 1.064483564774226, 1.07055262013105, 1.048166594536654, 1.057303911425504, 1.060089618360339, 1.067712419896691, 1.073761932664271, 1.058805404254182)
$data[19] = @(1.02, 1.051370310268492, 1.056889533120742, 1.064548882742808, 1.070617429368363, 1.04818453684334, 1.057346513585855, 1.060129202568869, 1.06776365483621, 1.073812751977241, 1.058848066914487)
$data[20] = @(1.02, 1.051059596858707, 1.056649219532449, 1.064256850417072, 1.070327669372361, 1.048104233893139, 1.057155993599898, 1.059952171127109, 1.067534553242658, 1.07358550614217, 1.058657276368297)
$data[21] = @(1.02, 1.050050503575995, 1.055868872652388, 1.06330892498608, 1.069387081870845, 1.047842097879084, 1.056536721078948, 1.059376602896158, 1.066790309360509, 1.072847242592532, 1.058037124409336)
$data[22] = @(1.02, 1.049416871482857, 1.055378963772218, 1.062714085727412, 1.068796818406927, 1.047676471045265, 1.056147463200223, 1.059014710003198, 1.066322831737228, 1.072383484629913, 1.057647313739775)
$data[23] = @(1.02, 1.049752693149122, 1.055638604332922, 1.063029310796794, 1.069109621336536, 1.047764349655453, 1.056353805952871, 1.059206556829593, 1.066570606555365, 1.072629291406426, 1.057853949522799)
$data[24] = @(1.02, 1.051076528314349, 1.056662314334797, 1.064272762026984, 1.070343457316458, 1.048108614830359, 1.057166377445026, 1.059961820337401, 1.067547038244907, 1.073597890196124, 1.058667674959676)
$data[25] = @(1.02, 1.052616306700964, 1.057853377142208, 1.065720703772174, 1.071780075021407, 1.048504558643713, 1.058109729265976, 1.060838174521405, 1.068682077311681, 1.074723662673717, 1.059612366448415)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}
